# Fruta / hortaliza, semanal
# The underlying data rows (2..28) got reshuffled: for each destination row,
# columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) are replaced with the
# values that used to live in a different (source) row. Columns A, B, C, E,
# F, G, H, I, N, O, Q, R are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (source row's current D/J/K/L/M/P values
# become the destination row's new D/J/K/L/M/P values)
$rowMap = @{
    2  = 16
    3  = 7
    4  = 11
    5  = 19
    6  = 13
    7  = 17
    8  = 15
    9  = 21
    11 = 20
    12 = 9
    13 = 14
    14 = 24
    15 = 4
    16 = 22
    17 = 5
    18 = 27
    19 = 2
    20 = 26
    21 = 6
    22 = 12
    23 = 18
    24 = 8
    26 = 28
    27 = 23
    28 = 3
}

# Columns involved in the reshuffle (1-based column index)
# D = 4, J = 10, K = 11, L = 12, M = 13, P = 16
$cols = @(4, 10, 11, 12, 13, 16)

# Snapshot every current value BEFORE any writes, since several rows source
# from each other (this is a set of permutation cycles, not simple pairwise
# swaps), so writes must never read from an already-overwritten row.
$snapshot = @{}
for ($r = 2; $r -le 28; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
